$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date strings in column A (slashes -> dashes) for rows 3 through 21.
# These are plain text cells (not real dates). Some of the new strings are
# ambiguous day/month values (day <= 12), which Excel's auto-detection would
# otherwise reinterpret as a date serial number, so those are entered with a
# leading apostrophe to force text, exactly as a user typing them in Excel
# would need to do.
$ws.Range("A3").Value = "28-07-2022"
$ws.Range("A4").Value = "'01-08-2022"
$ws.Range("A5").Value = "'04-08-2022"
$ws.Range("A6").Value = "'08-08-2022"
$ws.Range("A7").Value = "'11-08-2022"
$ws.Range("A8").Value = "15-08-2022"
$ws.Range("A9").Value = "18-08-2022"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("A13").Value = "'01-09-2022"
$ws.Range("A14").Value = "'05-09-2022"
$ws.Range("A15").Value = "'08-09-2022"
$ws.Range("A16").Value = "'12-09-2022"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A21").Value = "29-09-2022"

# Update attendance counts for row 4 (01-08-2022)
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("H4").Value = 0

# Update attendance counts for row 10 (22-08-2022)
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("H10").Value = 0

# Update attendance counts for row 12 (29-08-2022)
$ws.Range("D12").Value = 1
$ws.Range("G12").Value = 1
